# Generate Report for Handoff
# This script swaps the two file rows (9bff25ad.. <-> f1ddb85d..) on every
# sheet so the f1ddb85d item now sits in row 2 (still "In Translation"/"ht")
# and the 9bff25ad item moves to row 3 with a refreshed "Ready for
# handoff"/"mt" status and new handoff timestamps, matching a fresh
# handoff report run.

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/913bc1ca29d36529a276d990a6d5c17d72e623ae/e2e/"
$nameA = "9bff25ad-ad66-429a-b659-8b4f435ecf2a.md"
$nameB = "f1ddb85d-b6f1-496b-83f8-c73fdc7cea9b.md"

function Reset-Hyperlinks($ws) {
    $guard = 0
    while ($ws.Hyperlinks.Count -gt 0 -and $guard -lt 20) {
        foreach ($h in $ws.Hyperlinks) {
            $h.Delete()
        }
        $guard = $guard + 1
    }
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $nameB
$ws.Range("B2").Value = "e2e\" + $nameB
$ws.Range("A3").Value = $nameA
$ws.Range("B3").Value = "e2e\" + $nameA

$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-21 20:15:36"

Reset-Hyperlinks $ws
$null = $ws.Hyperlinks.Add($ws.Range("B2"), $urlBase + $nameA, "", "", "e2e\" + $nameB)
$null = $ws.Hyperlinks.Add($ws.Range("B3"), $urlBase + $nameB, "", "", "e2e\" + $nameA)

$ws.Columns.Item(5).ColumnWidth = 16.42
$ws.Columns.Item(6).ColumnWidth = 16.42

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $nameB
$ws.Range("G2").Value = "f1ddb85d-b6f1-496b-83f8-c73fdc7cea9b.f53450cc803b3f40a521345679de30cbd871d6ae.zh-cn.xlf"

$ws.Range("A3").Value = $nameA
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "9bff25ad-ad66-429a-b659-8b4f435ecf2a.d9c54e262a42a93f6175a9feba538ee36db94155.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-21 20:15:32"

Reset-Hyperlinks $ws
$null = $ws.Hyperlinks.Add($ws.Range("A2"), $urlBase + $nameA, "", "", $nameB)
$null = $ws.Hyperlinks.Add($ws.Range("A3"), $urlBase + $nameB, "", "", $nameA)

$ws.Columns.Item(3).ColumnWidth = 16.42

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $nameB
$ws.Range("G2").Value = "f1ddb85d-b6f1-496b-83f8-c73fdc7cea9b.f53450cc803b3f40a521345679de30cbd871d6ae.de-de.xlf"

$ws.Range("A3").Value = $nameA
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "9bff25ad-ad66-429a-b659-8b4f435ecf2a.d9c54e262a42a93f6175a9feba538ee36db94155.de-de.xlf"
$ws.Range("H3").Value = "2016-08-21 20:15:36"

Reset-Hyperlinks $ws
$null = $ws.Hyperlinks.Add($ws.Range("A2"), $urlBase + $nameA, "", "", $nameB)
$null = $ws.Hyperlinks.Add($ws.Range("A3"), $urlBase + $nameB, "", "", $nameA)

$ws.Columns.Item(3).ColumnWidth = 16.42
